$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.639.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.147.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.61%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.138.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000224"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.637.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.547.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.129.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "516.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "548.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0436"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.145.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.50%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0524"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.69%  "
